# week 11: update slides
# Applies:
#  1. Refresh cached "datetime1"/"datetimeFigureOut" field text (11/7/16 -> 4/3/17)
#     on the slide master, every slide layout, and the notes master.
#  2. Merge the two-run title on slide 1 ("Week 11: Listeners and Anonymous " + "Classes")
#     into a single run.
#  3. Split the "Has access to outer class’s static methods and fields" bullet on
#     slide 5 and append a new bullet ("Don’t need instance of outer class to access
#     inner class").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholders: walk every shape collection that can carry one
#    (ppPlaceholderDate = 16) and refresh its cached text.
# ---------------------------------------------------------------------------
$NEW_DATE = "4/3/17"

function Update-DatePlaceholders($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    $isDate = $false
    try {
      if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
    } catch {
      $isDate = $false
    }
    if ($isDate) {
      $sh.TextFrame.TextRange.Text = $NEW_DATE
    }
  }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
  $layout = $p.SlideMaster.CustomLayouts.Item($li)
  Update-DatePlaceholders $layout.Shapes
}

Update-DatePlaceholders $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2. Slide 1 subtitle: collapse the split title run into one run.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item("Subtitle 2")
$titleRange = $subtitle.TextFrame.TextRange
$titleChars = $titleRange.Characters(1, $titleRange.Length)
$titleChars.Text = "Week 11: Listeners and Anonymous Classes"

# ---------------------------------------------------------------------------
# 3. Slide 5 "Static Member Classes" content: split the 3rd bullet and add a
#    4th bullet.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$content = $slide5.Shapes.Item("Content Placeholder 1")
$bodyRange = $content.TextFrame.TextRange

$thirdPara = $bodyRange.Paragraphs(3, 1)
$splitText = "Has access to outer class’s static methods and "
$headChars = $bodyRange.Characters($thirdPara.Start, $splitText.Length)
$headChars.Text = $splitText

$newPara = "Don’t need instance of outer class to access inner class"
$null = $bodyRange.InsertAfter("`r" + $newPara)
